$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B6").Value = "ffffff"
$ws.Range("B6").Select()
